$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for existing rows 2-18
# from 45175 (2023-09-06) to 45177 (2023-09-08)
$ws.Range("C2:C18").Value = 45177

# Ensure row 18 carries an explicit row height (matches ht="15" customHeight="1")
$ws.Rows(18).RowHeight = 15

# Add the new data row 19: "A 41803-2023"
$ws.Range("A19").Value = "A 41803-2023"
$ws.Range("B19").Value = 45176
$ws.Range("C19").Value = 45177
$ws.Range("D19").Value = "VÄRMLANDS LÄN"
$ws.Range("E19").Value = "HAMMARÖ"
$ws.Range("F19").Value = "Övriga Aktiebolag"
$ws.Range("G19").Value = 1.2
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0

# Match the date number format used by columns B and C in other rows
$ws.Range("B19:C19").NumberFormat = "YYYY-MM-DD"

# Column R keeps the wrap-text style used throughout the table, with no value
$ws.Range("R19").WrapText = $true
